$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E1 value from 12 to 17
$ws.Range("E1").Value = 17

# Add new row 6 data
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "16-03-2024"
$ws.Range("D6").Value = "00:31:13"

# Add new row 7 data
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "16-03-2024"
$ws.Range("D7").Value = "01:57:58"
